# Quarterly indexing esoteric bug-fix operation
#
# For every data row (2..16) a new "Q(-1)" quarter-over-quarter error value
# needs to be inserted in front of the existing series (column B), pushing
# the previously recorded values one column to the right and dropping
# anything that would fall past column K (the series only keeps 10 columns,
# B:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to insert at column B for each row (2..16)
$newValues = @{
    2  = -1.025188112727922
    3  = 0.08364543516793629
    4  = -0.1538585523806955
    5  = 0.7495351060200912
    6  = 0.03849281619118239
    7  = -0.2590580299438133
    8  = 0.01855976243503714
    9  = 0.1467044301255134
    10 = -0.1819613811903656
    11 = 0.4718454808444464
    12 = -0.08594117411414147
    13 = -0.07695400962807622
    14 = -0.5068991247689255
    15 = 0.6215838649243215
    16 = -0.2766911554241067
}

$firstCol = 2   # column B
$lastCol  = 11  # column K

for ($r = 2; $r -le 16; $r++) {
    # Shift existing values one column to the right, starting from the
    # rightmost column so we never clobber a value before it has been read.
    for ($c = $lastCol; $c -ge ($firstCol + 1); $c--) {
        $srcVal = $ws.Cells.Item($r, $c - 1).Value2
        $ws.Cells.Item($r, $c).Value2 = $srcVal
    }
    # Write the new value into column B.
    $ws.Cells.Item($r, $firstCol).Value2 = $newValues[$r]
}
